$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.641.91'
$ws.Range('E2').Value = '  -0.45%  '
$ws.Range('D3').Value = '3.100.15'
$ws.Range('E3').Value = '  +1.35%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '523.30'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.50%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '141.39'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.099.18'
$ws.Range('E8').Value = '  +1.41%  '
$ws.Range('E9').Value = '  +0.59%  '
$ws.Range('E11').Value = '  +0.30%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.386'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +2.86%  '
$ws.Range('D13').Value = '3.633.83'
$ws.Range('E13').Value = '  +1.41%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.132'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.40%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '25.67'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.26%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000163'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.16%  '
$ws.Range('D17').Value = '57.733.52'
$ws.Range('E17').Value = '  -0.32%  '
$ws.Range('D18').Value = '3.103.61'
$ws.Range('E18').Value = '  +1.35%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.07'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.28%  '
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('E21').Value = '  -1.22%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '338.95'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.36%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('E24').Value = '  +2.27%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '66.65'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.86%  '
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').Value = '0.0₃0917'
$ws.Range('E28').Value = '  +1.67%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.47'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.18'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.80%  '
$ws.Range('E32').Value = '  +3.23%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '20.90'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.52%  '
$ws.Range('E34').Value = '  -0.37%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '155.66'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.72%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.61'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +1.90%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.12'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +2.09%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '27.08'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('E39').Value = '  -1.26%  '
$ws.Range('E40').Value = '  -2.29%  '
$ws.Range('E41').Value = '  +0.99%  '
$ws.Range('D42').Value = '3.138.94'
$ws.Range('E42').Value = '  +1.24%  '
$ws.Range('E43').Value = '  +10.53%  '
$ws.Range('E44').Value = '  +3.90%  '
$ws.Range('E45').Value = '  +0.42%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.00'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('D47').Value = '2.293.83'
$ws.Range('E47').Value = '  +0.44%  '
$ws.Range('E48').Value = '  +1.23%  '
$ws.Range('E49').Value = '  +4.93%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '20.52'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.73%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '6.02'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.65%  '
